# Insert a new data row at row 73 (pushing the existing rows 73-123 down to
# 74-124) on the single worksheet, then populate the new row with the
# "Macroferia Regional de Talca" / Cilantro record that was added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 73..123 down to 74..124, carrying over formatting (matches the
# blank-row-with-date-style behaviour Excel shows when inserting above a
# date-formatted cell).
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with the new record's values.
$ws.Cells.Item(73, 1).Value  = 5
$ws.Cells.Item(73, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(73, 3).Value  = "Maule"
$ws.Cells.Item(73, 4).Value  = 45176
$ws.Cells.Item(73, 5).Value  = 7
$ws.Cells.Item(73, 6).Value  = 100112040
$ws.Cells.Item(73, 7).Value  = "Cilantro"
$ws.Cells.Item(73, 8).Value  = "Sin especificar"
$ws.Cells.Item(73, 9).Value  = "Primera"
$ws.Cells.Item(73, 10).Value = 150
$ws.Cells.Item(73, 11).Value = 8000
$ws.Cells.Item(73, 12).Value = 8000
$ws.Cells.Item(73, 13).Value = 8000
$ws.Cells.Item(73, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(73, 15).Value = "Región Metropolitana"
$ws.Cells.Item(73, 16).Value = 222
$ws.Cells.Item(73, 17).Value = 36
$ws.Cells.Item(73, 18).Value = "Hortaliza"
